$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link / Volume label) -----------------------
# These are plain text cells; assigning a string keeps them as text.
$textUpdates = @{
    'B7' = 'KuCoinToken'
    'C7' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'E7' = '6KuCoinTokenKCS'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'E8' = '7MXTokenMX'
    'B9' = 'FTXToken'
    'C9' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'E9' = '8FTXTokenFTT'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'E10' = '9WazirXWRX'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'E11' = '10MandalaExchangeTokenMDX'
    'B12' = 'LiechtensteinCryptoassetsExchange'
    'C12' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'E12' = '11LiechtensteinCryptoassetsExchangeLCX'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'E13' = '12BitrueCoinBTR'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'E14' = '13BitMartTokenBMX'
    'B15' = 'MCDex'
    'C15' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'E15' = '14MCDexMCB'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'E16' = '15BitForexTokenBF'
    'B17' = 'CoinExToken'
    'C17' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'E17' = '16CoinExTokenCET'
    'B18' = 'One'
    'C18' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'E18' = '17OneONEWorstin24h'
    'B19' = 'TigerCash'
    'C19' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'E19' = '18TigerCashTCH'
    'B20' = 'HotbitToken'
    'C20' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'E20' = '19HotbitTokenHTB'
    'B21' = 'BitKan'
    'C21' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'E21' = '20BitKanKAN'
    'B22' = 'NitroEx'
    'C22' = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
    'E22' = '21NitroExNTX'
    'B23' = 'LEO'
    'C23' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'E23' = '22LEOLEO'
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# --- Price column (D) --------------------------------------------------
# Prices are stored as text in this sheet (e.g. "0.0001500"), so force the
# number format to Text before writing the value to preserve the exact
# string representation (including trailing zeros) instead of letting Excel
# coerce it to a numeric value.
$priceUpdates = @{
    'D2' = '247.24'
    'D3' = '22.62'
    'D4' = '5.286'
    'D5' = '0.05724'
    'D6' = '3.425'
    'D7' = '6.335'
    'D8' = '0.8067'
    'D9' = '0.8636'
    'D10' = '0.1422'
    'D11' = '0.07341'
    'D12' = '0.03043'
    'D13' = '0.03112'
    'D14' = '0.09397'
    'D15' = '3.880'
    'D16' = '0.001579'
    'D17' = '0.04812'
    'D18' = '0.0005850'
    'D19' = '0.006153'
    'D20' = '0.005165'
    'D21' = '0.0009969'
    'D22' = '0.0001500'
    'D23' = '3.735'
    'D26' = '0.1296'
    'D40' = '0.03940'
    'D41' = '0.006757'
    'D42' = '0.1069'
    'D44' = '0.007973'
    'D45' = '0.00005608'
    'D47' = '0.3600'
    'D48' = '0.1817'
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}
